$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.000000008427485376216736767674
$ws.Range("C2").Value = 0.004309184025731883
$ws.Range("D2").Value = 2938.103010863317
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("G2").Value = 74455.9988974649

# Row 3
$ws.Range("B3").Value = 0.0003714022599530242
$ws.Range("C3").Value = 0.000002220651329265522090883030
$ws.Range("D3").Value = 16.98373111632243
$ws.Range("E3").Value = 71517.89157740913
$ws.Range("G3").Value = 71534.87568214837

# Row 4
$ws.Range("B4").Value = 0.3464964993005633
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 16.98373111632243
$ws.Range("E4").Value = 6.48142807727062
$ws.Range("G4").Value = 25.46489215179242

# Row 5
$ws.Range("B5").Value = 0.7287194209349384
$ws.Range("C5").Value = 9.226618575922256
$ws.Range("D5").Value = 3.082599426703578
$ws.Range("E5").Value = 6.48142807727062
$ws.Range("G5").Value = 19.51936550083139

